$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 34-63: only columns D,L,M,N,O,P,R,S change (data shifted down by 2,
#     rows 34-35 replaced with two brand new weekly records) ---

# Row 34
$ws.Range("D34").Value = 44566
$ws.Range("L34").Value = "Especial"
$ws.Range("M34").Value = 410
$ws.Range("N34").Value = 8000
$ws.Range("O34").Value = 8000
$ws.Range("P34").Value = 8000
$ws.Range("R34").Value = "Provincia de Curicó"
$ws.Range("S34").Value = 4000

# Row 35
$ws.Range("D35").Value = 44566
$ws.Range("L35").Value = "Primera"
$ws.Range("M35").Value = 450
$ws.Range("N35").Value = 7000
$ws.Range("O35").Value = 7000
$ws.Range("P35").Value = 7000
$ws.Range("R35").Value = "Provincia de Curicó"
$ws.Range("S35").Value = 3500

# Row 36
$ws.Range("D36").Value = 44273
$ws.Range("L36").Value = "Primera"
$ws.Range("M36").Value = 210
$ws.Range("N36").Value = 6000
$ws.Range("O36").Value = 6000
$ws.Range("P36").Value = 6000
$ws.Range("R36").Value = "Provincia de Linares"
$ws.Range("S36").Value = 3000

# Row 37
$ws.Range("D37").Value = 44364
$ws.Range("L37").Value = "Primera"
$ws.Range("M37").Value = 75
$ws.Range("N37").Value = 10000
$ws.Range("O37").Value = 10000
$ws.Range("P37").Value = 10000
$ws.Range("R37").Value = "Provincia de Curicó"
$ws.Range("S37").Value = 5000

# Row 38
$ws.Range("D38").Value = 44280
$ws.Range("L38").Value = "Primera"
$ws.Range("M38").Value = 260
$ws.Range("N38").Value = 8000
$ws.Range("O38").Value = 8000
$ws.Range("P38").Value = 8000
$ws.Range("R38").Value = "Provincia de Linares"
$ws.Range("S38").Value = 4000

# Row 39
$ws.Range("D39").Value = 44209
$ws.Range("L39").Value = "Primera"
$ws.Range("M39").Value = 370
$ws.Range("N39").Value = 5800
$ws.Range("O39").Value = 6000
$ws.Range("P39").Value = 5935
$ws.Range("R39").Value = "Provincia de Linares"
$ws.Range("S39").Value = 2968

# Row 40
$ws.Range("D40").Value = 44525
$ws.Range("L40").Value = "Primera"
$ws.Range("M40").Value = 350
$ws.Range("N40").Value = 8000
$ws.Range("O40").Value = 8000
$ws.Range("P40").Value = 8000
$ws.Range("R40").Value = "Provincia de Curicó"
$ws.Range("S40").Value = 4000

# Row 41
$ws.Range("D41").Value = 44165
$ws.Range("L41").Value = "Especial"
$ws.Range("M41").Value = 25
$ws.Range("N41").Value = 10000
$ws.Range("O41").Value = 10000
$ws.Range("P41").Value = 10000
$ws.Range("R41").Value = "Provincia de Linares"
$ws.Range("S41").Value = 5000

# Row 42
$ws.Range("D42").Value = 44306
$ws.Range("L42").Value = "Primera"
$ws.Range("M42").Value = 200
$ws.Range("N42").Value = 7000
$ws.Range("O42").Value = 7000
$ws.Range("P42").Value = 7000
$ws.Range("R42").Value = "Provincia de Curicó"
$ws.Range("S42").Value = 3500

# Row 43
$ws.Range("D43").Value = 44533
$ws.Range("L43").Value = "Primera"
$ws.Range("M43").Value = 350
$ws.Range("N43").Value = 10000
$ws.Range("O43").Value = 10000
$ws.Range("P43").Value = 10000
$ws.Range("R43").Value = "Provincia de Curicó"
$ws.Range("S43").Value = 5000

# Row 44
$ws.Range("D44").Value = 44539
$ws.Range("L44").Value = "Primera"
$ws.Range("M44").Value = 450
$ws.Range("N44").Value = 6000
$ws.Range("O44").Value = 6000
$ws.Range("P44").Value = 6000
$ws.Range("R44").Value = "Provincia de Curicó"
$ws.Range("S44").Value = 3000

# Row 45
$ws.Range("D45").Value = 44210
$ws.Range("L45").Value = "Primera"
$ws.Range("M45").Value = 400
$ws.Range("N45").Value = 5800
$ws.Range("O45").Value = 6000
$ws.Range("P45").Value = 5910
$ws.Range("R45").Value = "Provincia de Linares"
$ws.Range("S45").Value = 2955

# Row 46
$ws.Range("D46").Value = 44559
$ws.Range("L46").Value = "Primera"
$ws.Range("M46").Value = 450
$ws.Range("N46").Value = 8000
$ws.Range("O46").Value = 8000
$ws.Range("P46").Value = 8000
$ws.Range("R46").Value = "Provincia de Curicó"
$ws.Range("S46").Value = 4000

# Row 47
$ws.Range("D47").Value = 44286
$ws.Range("L47").Value = "Primera"
$ws.Range("M47").Value = 100
$ws.Range("N47").Value = 8000
$ws.Range("O47").Value = 8000
$ws.Range("P47").Value = 8000
$ws.Range("R47").Value = "Provincia de Linares"
$ws.Range("S47").Value = 4000

# Row 48
$ws.Range("D48").Value = 44279
$ws.Range("L48").Value = "Primera"
$ws.Range("M48").Value = 150
$ws.Range("N48").Value = 8000
$ws.Range("O48").Value = 8000
$ws.Range("P48").Value = 8000
$ws.Range("R48").Value = "Provincia de Linares"
$ws.Range("S48").Value = 4000

# Row 49
$ws.Range("D49").Value = 44208
$ws.Range("L49").Value = "Primera"
$ws.Range("M49").Value = 300
$ws.Range("N49").Value = 6000
$ws.Range("O49").Value = 6000
$ws.Range("P49").Value = 6000
$ws.Range("R49").Value = "Provincia de Linares"
$ws.Range("S49").Value = 3000

# Row 50
$ws.Range("D50").Value = 44264
$ws.Range("L50").Value = "Primera"
$ws.Range("M50").Value = 220
$ws.Range("N50").Value = 6000
$ws.Range("O50").Value = 6000
$ws.Range("P50").Value = 6000
$ws.Range("R50").Value = "Región de O'Higgins"
$ws.Range("S50").Value = 3000

# Row 51
$ws.Range("D51").Value = 44281
$ws.Range("L51").Value = "Primera"
$ws.Range("M51").Value = 220
$ws.Range("N51").Value = 7500
$ws.Range("O51").Value = 7500
$ws.Range("P51").Value = 7500
$ws.Range("R51").Value = "Provincia de Linares"
$ws.Range("S51").Value = 3750

# Row 52
$ws.Range("D52").Value = 44187
$ws.Range("L52").Value = "Primera"
$ws.Range("M52").Value = 220
$ws.Range("N52").Value = 7000
$ws.Range("O52").Value = 7000
$ws.Range("P52").Value = 7000
$ws.Range("R52").Value = "Provincia de Linares"
$ws.Range("S52").Value = 3500

# Row 53
$ws.Range("D53").Value = 44187
$ws.Range("L53").Value = "Segunda"
$ws.Range("M53").Value = 260
$ws.Range("N53").Value = 5000
$ws.Range("O53").Value = 5000
$ws.Range("P53").Value = 5000
$ws.Range("R53").Value = "Provincia de Linares"
$ws.Range("S53").Value = 2500

# Row 54
$ws.Range("D54").Value = 44238
$ws.Range("L54").Value = "Primera"
$ws.Range("M54").Value = 200
$ws.Range("N54").Value = 6000
$ws.Range("O54").Value = 6000
$ws.Range("P54").Value = 6000
$ws.Range("R54").Value = "Provincia de Curicó"
$ws.Range("S54").Value = 3000

# Row 55
$ws.Range("D55").Value = 44238
$ws.Range("L55").Value = "Primera"
$ws.Range("M55").Value = 150
$ws.Range("N55").Value = 6000
$ws.Range("O55").Value = 6000
$ws.Range("P55").Value = 6000
$ws.Range("R55").Value = "Provincia de Linares"
$ws.Range("S55").Value = 3000

# Row 56
$ws.Range("D56").Value = 44196
$ws.Range("L56").Value = "Primera"
$ws.Range("M56").Value = 550
$ws.Range("N56").Value = 6500
$ws.Range("O56").Value = 7000
$ws.Range("P56").Value = 6818
$ws.Range("R56").Value = "Provincia de Linares"
$ws.Range("S56").Value = 3409

# Row 57
$ws.Range("D57").Value = 44188
$ws.Range("L57").Value = "Primera"
$ws.Range("M57").Value = 260
$ws.Range("N57").Value = 6500
$ws.Range("O57").Value = 7000
$ws.Range("P57").Value = 6692
$ws.Range("R57").Value = "Provincia de Linares"
$ws.Range("S57").Value = 3346

# Row 58
$ws.Range("D58").Value = 44188
$ws.Range("L58").Value = "Segunda"
$ws.Range("M58").Value = 340
$ws.Range("N58").Value = 5000
$ws.Range("O58").Value = 5500
$ws.Range("P58").Value = 5206
$ws.Range("R58").Value = "Provincia de Linares"
$ws.Range("S58").Value = 2603

# Row 59
$ws.Range("D59").Value = 44224
$ws.Range("L59").Value = "Primera"
$ws.Range("M59").Value = 420
$ws.Range("N59").Value = 6500
$ws.Range("O59").Value = 7000
$ws.Range("P59").Value = 6786
$ws.Range("R59").Value = "Provincia de Linares"
$ws.Range("S59").Value = 3393

# Row 60
$ws.Range("D60").Value = 44195
$ws.Range("L60").Value = "Primera"
$ws.Range("M60").Value = 408
$ws.Range("N60").Value = 6509
$ws.Range("O60").Value = 7000
$ws.Range("P60").Value = 6774
$ws.Range("R60").Value = "Provincia de Linares"
$ws.Range("S60").Value = 3387

# Row 61
$ws.Range("D61").Value = 44302
$ws.Range("L61").Value = "Primera"
$ws.Range("M61").Value = 150
$ws.Range("N61").Value = 7000
$ws.Range("O61").Value = 7000
$ws.Range("P61").Value = 7000
$ws.Range("R61").Value = "Provincia de Curicó"
$ws.Range("S61").Value = 3500

# Row 62
$ws.Range("D62").Value = 44239
$ws.Range("L62").Value = "Primera"
$ws.Range("M62").Value = 150
$ws.Range("N62").Value = 6000
$ws.Range("O62").Value = 6000
$ws.Range("P62").Value = 6000
$ws.Range("R62").Value = "Provincia de Curicó"
$ws.Range("S62").Value = 3000

# Row 63
$ws.Range("D63").Value = 44357
$ws.Range("L63").Value = "Primera"
$ws.Range("M63").Value = 100
$ws.Range("N63").Value = 10000
$ws.Range("O63").Value = 10000
$ws.Range("P63").Value = 10000
$ws.Range("R63").Value = "Provincia de Curicó"
$ws.Range("S63").Value = 5000

# --- Append two brand new rows (64, 65) that did not exist before ---

# Row 64
$ws.Range("A64").Value = 9
$ws.Range("B64").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C64").Value = "Metropolitana"
$ws.Range("D64").Value = 44365
$ws.Range("E64").Value = 13
$ws.Range("F64").Value = "Fruta"
$ws.Range("G64").Value = 100101
$ws.Range("H64").Value = "Berries"
$ws.Range("I64").Value = 100101004
$ws.Range("J64").Value = "Frambuesa"
$ws.Range("K64").Value = "Sin especificar"
$ws.Range("L64").Value = "Primera"
$ws.Range("M64").Value = 50
$ws.Range("N64").Value = 10000
$ws.Range("O64").Value = 10000
$ws.Range("P64").Value = 10000
$ws.Range("Q64").Value = "`$/bandeja 2 kilos"
$ws.Range("R64").Value = "Provincia de Curicó"
$ws.Range("S64").Value = 5000
$ws.Range("T64").Value = 2
$ws.Range("D64").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Row 65
$ws.Range("A65").Value = 9
$ws.Range("B65").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C65").Value = "Metropolitana"
$ws.Range("D65").Value = 44194
$ws.Range("E65").Value = 13
$ws.Range("F65").Value = "Fruta"
$ws.Range("G65").Value = 100101
$ws.Range("H65").Value = "Berries"
$ws.Range("I65").Value = 100101004
$ws.Range("J65").Value = "Frambuesa"
$ws.Range("K65").Value = "Sin especificar"
$ws.Range("L65").Value = "Primera"
$ws.Range("M65").Value = 190
$ws.Range("N65").Value = 5800
$ws.Range("O65").Value = 6000
$ws.Range("P65").Value = 5916
$ws.Range("Q65").Value = "`$/bandeja 2 kilos"
$ws.Range("R65").Value = "Provincia de Linares"
$ws.Range("S65").Value = 2958
$ws.Range("T65").Value = 2
$ws.Range("D65").NumberFormat = "YYYY-MM-DD HH:MM:SS"
